$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("steps")

# Insert a new blank row at row 62 (this pushes existing row 63.. content
# down by one row each, matching the old row N -> new row N+1 mapping
# seen through the rest of the sheet).
$ws.Rows.Item(62).Insert()

# The second new line lands in a row that was already blank (the gap
# between "git commit" and "git push heroku master" shifted down by the
# insert above), so no further row insertion is required - just set the
# value directly. Set this one first so it becomes the earlier shared
# string entry.
$ws.Range("B79").Value = "heroku git:remote -a bhver2"

# Fill in the newly inserted row with the other new line, styled in red
# text (same red Calibri font used elsewhere in the sheet).
$ws.Range("B62").Value = "use app2.py to run locally , app.py runs on heroku"
$ws.Range("B62").Font.Color = 255

# Restore the view state recorded after the edit.
$ws.Application.ActiveWindow.ScrollRow = 53
$ws.Range("B62").Select()
